# Season record: add Wins / Losses / Ties columns (AD:AF) to the sheet.
# The team finished the season 72-90-0, so every player row gets the same
# season record values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 48

# --- Header row (row 1) -----------------------------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting of the existing header cells (bold, centered,
# bordered) by copying the format from the last existing header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (rows 2-48) ---------------------------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 72
    $ws.Cells.Item($r, 31).Value = 90
    $ws.Cells.Item($r, 32).Value = 0
}
